$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The date/time values in B7:B9 were being auto-converted into Excel serial
# date numbers. Fix the "Excel shenanigans" by forcing the cells back to a
# plain General number format and writing the literal timestamp text.
$ws.Range("B7:B9").NumberFormat = "general"

$ws.Range("B7").Value = "28/08/2017 07:08"
$ws.Range("B8").Value = "28/08/2017 07:30"
$ws.Range("B9").Value = "28/08/2017 07:33"

$ws.Range("B7").Select()
